# Trading update: 2026-02-17 08:20:59
# Appends a new open trade (row 20) to both the "All Trades" and
# "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 20

    $ws.Cells.Item($r, 1).Value = 19

    # Use a leading apostrophe so Excel stores the ISO-looking date string as
    # plain text instead of auto-converting it to a date serial number, then
    # reset the style back to Normal so no quote-prefix styling lingers.
    $ws.Cells.Item($r, 2).Value = "'2026-02-17"
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value = "08:19:37"
    $ws.Cells.Item($r, 4).Value = "MarketMaking"
    $ws.Cells.Item($r, 5).Value = "DOWN"
    $ws.Cells.Item($r, 6).Value = 0.98
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = "OPEN"
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 99.99697504264921
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0.6
    $ws.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($r, 16).Value = ""
    $ws.Cells.Item($r, 17).Value = 0
}
